$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "`"GraphicRasm5`""
$ws.Range("C3").Value = "LoyihaIshiBir"

$ws.Range("D3").Select()
